# Update "想去人数" (number of people interested) figures that changed
# between the previous data scrape and the latest one (commit 456a3b4).
#
# Sheet 1 = 展览 (Exhibitions), Sheet 4 = 全部类型 (All types) both contain
# the same events (全部类型 aggregates everything) so both need updating.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value = 7930
$ws1.Range("F6").Value = 8674
$ws1.Range("F10").Value = 5803
$ws1.Range("F12").Value = 2804
$ws1.Range("F17").Value = 636
$ws1.Range("F19").Value = 3995
$ws1.Range("F20").Value = 3995
$ws1.Range("F27").Value = 5759
$ws1.Range("F28").Value = 5759
$ws1.Range("F34").Value = 3294
$ws1.Range("F42").Value = 3696

# ---- Sheet 4: 全部类型 ----
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 7930
$ws4.Range("F6").Value = 8674
$ws4.Range("F10").Value = 5803
$ws4.Range("F12").Value = 2804
$ws4.Range("F18").Value = 636
$ws4.Range("F20").Value = 3995
$ws4.Range("F21").Value = 3995
$ws4.Range("F28").Value = 5759
$ws4.Range("F29").Value = 5759
$ws4.Range("F35").Value = 3294
$ws4.Range("F44").Value = 3696
